$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking") - correct the Right / Wrong marks
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total") - correct the total marks and the max-score display text
$ws.Range("B12").Value = 48
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "46 / 112"
